# Nanobody_Library_ConditionTest_1 - bug fix + selection/view touch-up
#
# The underlying growth-curve data had two "bad" wells (columns G11/G12,
# i.e. spreadsheet columns CG/CH, and H11/H12, i.e. columns CS/CT) whose
# sensor readings had drifted. The fix replaces those four columns'
# values, row by row, with the value from the neighboring good well in
# column G10 (spreadsheet column CF) for every data row (rows 2-42).
#
# Afterwards the view is left scrolled/selected over the fixed columns
# (H11:H12, i.e. CT2:CT42) the way the author would have been reviewing
# the fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers (1-based): CF=84, CG=85, CH=86, CS=97, CT=98
$srcCol = 84
$badCols = @(85, 86, 97, 98)

for ($r = 2; $r -le 42; $r++) {
    $good = $ws.Cells.Item($r, $srcCol).Value()
    foreach ($c in $badCols) {
        $ws.Cells.Item($r, $c).Value = $good
    }
}

# Leave the window scrolled toward the repaired columns and the last
# repaired column (CT) selected, matching the reviewed state.
$excel.ActiveWindow.ScrollColumn = 60
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("CT2:CT42").Select()
